# Fix Training Data Issue
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF holds the date string for each row; convert from "4-30-2012-13" to "2013-04-30".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value2 -eq "4-30-2012-13") {
        # The leading apostrophe is Excel's standard "force text" quote
        # prefix: it keeps the replacement as the literal string
        # "2013-04-30" instead of being auto-parsed into a date serial
        # value (which is what a bare "2013-04-30" assignment would do).
        $cell.Formula = "'2013-04-30"
    }
}
